# Fix the F column labels for rows 78-81: they should all read
# "New CRM opened 12/11/2019" (same text as F77), not the accidental
# "2020/2021/2022/2023" variants that got typed in by mistake.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$label = "New CRM opened 12/11/2019"
$ws.Range("F78").Value = $label
$ws.Range("F79").Value = $label
$ws.Range("F80").Value = $label
$ws.Range("F81").Value = $label

# Append two more rows of sample data (rows 82 and 83), continuing the
# date/CRM-value/batch-value series below row 81.

# Copy row 81's date formatting down to A82/A83 first (copy/paste-format
# keeps using the workbook's existing date style instead of minting a new
# one), then fill in the actual values for every column.
$ws.Range("A81").Copy()
$ws.Range("A82").PasteSpecial(-4122)
$ws.Range("A83").PasteSpecial(-4122)

$ws.Range("A82").Value = 43816
$ws.Range("B82").Value = 2203.8582110000002
$ws.Range("C82").Value = 2207.0300000000002
$ws.Range("E82").Value = 169
$ws.Range("F82").Value = $label

$ws.Range("A83").Value = 43817
$ws.Range("B83").Value = 2208.9061499999998
$ws.Range("C83").Value = 2207.0300000000002
$ws.Range("E83").Value = 169
$ws.Range("F83").Value = $label

# Fill column D (the "% off" formula) for the two new rows in one shot so
# they pick up the same relative formula as the rest of the column.
$ws.Range("D82:D83").Formula = "=100*(B82-C82)/C82"

$excel.CutCopyMode = 0

# Update selection to match final state (cell B83 selected)
$ws.Range("B83").Select()
